$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($sheet, $row, $vals)
    $n = $vals.Length
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $vals[$i]
    }
    $startCol = 4  # column D
    $endCol = $startCol + $n - 1
    $rng = $sheet.Range($sheet.Cells.Item($row, $startCol), $sheet.Cells.Item($row, $endCol))
    $rng.Value = $arr
}

# Insert a new column before D. This shifts old D:K -> E:L and keeps
# their number formatting / styles intact.
$ws.Columns("D:D").Insert()

# The freshly inserted column D does not inherit the numeric / date
# style used by the rest of the table (it copies the sheet's generic
# column style). Re-apply the same formatting as the (now shifted)
# column E so D matches the look of the table.
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Set-RowValues $ws 7 @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908)
Set-RowValues $ws 8 @(33467800, 32869800, 31088100, 27806300, 25624000, 26380200, 28279800, 27926100)
Set-RowValues $ws 9 @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
Set-RowValues $ws 10 @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
Set-RowValues $ws 12 @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
Set-RowValues $ws 13 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 14 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 15 @(-1354200, -1557300, -1598800, -1427200, -1285800, -1228600, -1097300, -950800)
Set-RowValues $ws 17 @(18197600, 17070000, 15983900, 14597100, 14357000, 17081200, 20420200, 17243300)
Set-RowValues $ws 18 @(15270300, 15799900, 15104200, 13209200, 11267000, 9299100, 7859500, 10682800)
Set-RowValues $ws 20 @(-5796200, -8023400, -7932500, -8044700, -6801500, -8228700, -6084600, -6694200)
Set-RowValues $ws 21 @(10832700, 9336500, 8775600, 6595100, 5753300, 2306400, 13471800, 12367100)
Set-RowValues $ws 22 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 23 @(9474100, 7776500, 7171800, 5164500, 4465500, 1070400, 1775000, 3988600)
Set-RowValues $ws 24 @(2575000, 2433600, 1906300, 1429400, 1007500, -18000, -394900, 185500)
Set-RowValues $ws 25 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 26 @(6899100, 5342900, 5265500, 3735100, 3458000, 1088300, 2169900, 3803100)
Set-RowValues $ws 27 @(5620000, 3610600, 3607200, 2726400, 2937400, 244600, 1546100, 3283100)
Set-RowValues $ws 28 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 29 @("NA", "NA", "NA", 0, 0, 2093600, 524000, 343900)
Set-RowValues $ws 30 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 31 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 32 @(5796200, 8023400, 7932500, 8044700, 6801500, 8228700, 6084600, 6694200)
Set-RowValues $ws 33 @(5620000, 3610600, 3607200, 2726400, 2937400, 2338200, 2070100, 3627100)
Set-RowValues $ws 34 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 35 @(5620000, 3610600, 3607200, 2726400, 2937400, 2338200, 2070100, 3627100)
Set-RowValues $ws 38 @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908)
Set-RowValues $ws 41 @(79999000, 85541600, 90102500, 101924000, 49421400, 50252800, 9238500, 10566600)
Set-RowValues $ws 42 @(174799000, 156365000, 188009000, 249605000, 231105000, 194074000, 50993300, 45849000)
Set-RowValues $ws 43 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 44 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 45 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 46 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 47 @(1770500, 1781700, 858300, 986200, 5059100, 5320500, 12097300, 10915300)
Set-RowValues $ws 48 @(8110900, 8068200, 10032800, 11157100, 8774000, 8453100, 8495700, 8365700)
Set-RowValues $ws 49 @(9328200, 9496500, 10979800, 11278200, 8270200, 7583500, 8002000, 8075800)
Set-RowValues $ws 50 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 51 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 52 @(19429500, 43284100, 22433100, 21594900, 11658600, 10324600, "NA", "NA")
Set-RowValues $ws 53 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 54 @(759238000, 774239000, 821135000, 841330000, 709033000, 653780000, 696837000, 684141000)
Set-RowValues $ws 57 @(6337000, 6241600, 7326600, 7827000, 3753100, 2396600, 8515900, 8697900)
Set-RowValues $ws 58 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 59 @(1380000, 1249900, 2020700, 1960100, 1099600, 1336300, 6868800, 5113100)
Set-RowValues $ws 60 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 61 @(71773700, 71712000, 85692000, 91980700, 80997600, 83785700, 100286000, 126617000)
Set-RowValues $ws 62 @(9893700, 10839500, 13983400, 13766800, 11916700, 9413500, 6481700, 6546300)
Set-RowValues $ws 63 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 64 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 65 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 66 @(706381000, 722242000, 767993000, 788271000, 653944000, 606439000, 650353000, 639341000)
Set-RowValues $ws 68 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 69 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 70 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 71 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 72 @(31734400, 30401400, 30401400, 28197900, 26428500, 24516600, 21571400, 21373900)
Set-RowValues $ws 73 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 74 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 75 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 76 @(52856900, 51997500, 53141900, 53058900, 55088600, 47341200, 46484000, 44799600)
Set-RowValues $ws 77 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 80 @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908)
Set-RowValues $ws 81 @(5620000, 3610600, 3607200, 2726400, 2937400, 2338200, 2070100, 3627100)
Set-RowValues $ws 83 @(1355400, 1556200, 1600000, 1427200, 1284700, 1233100, 11668700, 8372800)
Set-RowValues $ws 84 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 85 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 86 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 87 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 88 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 89 @(9720900, 2306800, 7433200, 25919100, -6942900, -561000, 10914700, 20168400)
Set-RowValues $ws 91 @(-1058000, -871800, -1472100, -2435800, -1592100, -1404700, -1890600, -2244300)
Set-RowValues $ws 92 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 93 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 94 @(8432900, 3256000, -628300, -4949100, -1291400, 3389500, -1189300, -5977000)
Set-RowValues $ws 96 @(-2364000, -1905100, -1794100, -986200, -926800, -1430500, -1423800, -1210200)
Set-RowValues $ws 97 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 98 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 99 @(0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 100 @(-5713200, -110000, -1248800, 142500, 3542100, -1487800, -3918000, -1489600)
Set-RowValues $ws 101 @(-2802700, -4786400, -3885500, -7608200, 813400, -2001600, 528500, -1125700)
Set-RowValues $ws 102 @(9637900, 666500, 1670600, 13504300, -3878700, -660900, 6335900, 11576100)
